$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 600
$ws1.Range("F8").Value = 1534
$ws1.Range("F9").Value = 151
$ws1.Range("F10").Value = 151
$ws1.Range("F11").Value = 1440
$ws1.Range("F13").Value = 594
$ws1.Range("F14").Value = 1741
$ws1.Range("F15").Value = 1791
$ws1.Range("F16").Value = 840
$ws1.Range("F21").Value = 4
$ws1.Range("F22").Value = 1193
$ws1.Range("F24").Value = 438
$ws1.Range("F25").Value = 83
$ws1.Range("F26").Value = 4721
$ws1.Range("F30").Value = 64
$ws1.Range("F31").Value = 109

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 50
$ws2.Range("F5").Value = 24
$ws2.Range("F7").Value = 62

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 36

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 36
$ws4.Range("F4").Value = 50
$ws4.Range("F8").Value = 24
$ws4.Range("F10").Value = 62
$ws4.Range("F15").Value = 600
$ws4.Range("F16").Value = 1534
$ws4.Range("F17").Value = 151
$ws4.Range("F18").Value = 151
$ws4.Range("F20").Value = 1440
$ws4.Range("F22").Value = 594
$ws4.Range("F23").Value = 1741
$ws4.Range("F24").Value = 1791
$ws4.Range("F25").Value = 840
$ws4.Range("F31").Value = 4
$ws4.Range("F33").Value = 1193
$ws4.Range("F35").Value = 438
$ws4.Range("F36").Value = 83
$ws4.Range("F37").Value = 4721
$ws4.Range("F43").Value = 64
$ws4.Range("F44").Value = 109
